$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old content in column B (rows 1-7) so those shared strings are pruned
$ws.Range("B1:B7").ClearContents()

# Write the new cell values in reading order so the regenerated shared-string
# table comes out in the same order as the target workbook
$ws.Range("A1").Value = "This is Lesson 0. If you don't want a lesson 0, fair enough! Just start in row 2"
$ws.Range("A2").Value = "This is Lesson 1. Lesson 1 should be in row 2 etc."
$ws.Range("B2").Value = "If your lessons are over here, set Column Number to 1 in lesson_indexes.json. If they are in another column, set the appropriate column number"
$ws.Range("A3").Value = "Lesson 2"
$ws.Range("A4").Value = "Lesson 3"
$ws.Range("A5").Value = "Lesson 4"
$ws.Range("A6").Value = "You get the idea"
$ws.Range("A7").Value = "Lorum Impusm"
$ws.Range("A8").Value = "Dolor sit amet"

# Apply cell formatting -- B2 first so it claims cellXfs index 1 (center + wrap),
# then A1 so it claims cellXfs index 2 (wrap only), matching the target style order
$ws.Range("B2").HorizontalAlignment = -4108
$ws.Range("B2").WrapText = $true

$ws.Range("A1").WrapText = $true

# Row heights
$ws.Rows.Item(1).RowHeight = 30.75
$ws.Rows.Item(2).RowHeight = 43.5

# Column widths (closest attainable values in this engine's quantized model)
$ws.Columns.Item(1).ColumnWidth = 46.333333
$ws.Columns.Item(2).ColumnWidth = 43.5

# Selection matches the target workbook
$ws.Range("B2").Select() | Out-Null
